$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "KOORUI Monitor de juegos FHD de 24.5 pulgadas, monitores de computadora Full HD (1920 x 1080), VA, 1ms, compatible con FreeSync y G-Sync, 2 HDMI (144 Hz) y DisplayPort (170 Hz), sRGB 99%, VESA, inclinación ajustable, cuidado de los ojos"
